$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 5666.25
$ws.Range("I8").Value = 8062.3335
$ws.Range("K8").Value = 24187.0005
$ws.Range("M8").Value = -24048.0005
$ws.Range("H18").Value = 3582.3333
$ws.Range("I18").Value = 2698.8
$ws.Range("K18").Value = 2698.8
$ws.Range("M18").Value = -2414.8
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H86").Value = 4885.4546
$ws.Range("I86").Value = 4017
$ws.Range("J86").Value = 5927.6
$ws.Range("K86").Value = 4017
$ws.Range("L86").Value = 5927.6
$ws.Range("M86").Value = -2894
$ws.Range("N86").Value = -8173.6
$ws.Range("H89").Value = 4885.4546
$ws.Range("I89").Value = 4017
$ws.Range("J89").Value = 5927.6
$ws.Range("K89").Value = 20085
$ws.Range("L89").Value = 29638
$ws.Range("M89").Value = -14469
$ws.Range("N89").Value = -40870
$ws.Range("H98").Value = 7478.636
$ws.Range("J98").Value = 2238
$ws.Range("L98").Value = 2238
$ws.Range("N98").Value = -5234
$ws.Range("H122").Value = 7478.636
$ws.Range("J122").Value = 2238
$ws.Range("L122").Value = 6714
$ws.Range("N122").Value = -11614
$ws.Range("H137").Value = 1851.3636
$ws.Range("I137").Value = 1933.1578
$ws.Range("K137").Value = 5799.4734
$ws.Range("M137").Value = -3249.4734

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10192.0625
$ws.Range("I32").Value = 10192.0625
$ws.Range("K32").Value = 10192.0625
$ws.Range("M32").Value = -9905.0625

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4887
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 335.27274
$ws.Range("I22").Value = 348.66666
$ws.Range("K22").Value = 348.66666
$ws.Range("M22").Value = -175.66666
$ws.Range("H94").Value = 4063
$ws.Range("J94").Value = 1300
$ws.Range("L94").Value = 1300
$ws.Range("N94").Value = -2202
$ws.Range("H96").Value = 20717.3
$ws.Range("I96").Value = 13575.333
$ws.Range("J96").Value = 84995
$ws.Range("K96").Value = 13575.333
$ws.Range("L96").Value = 84995
$ws.Range("M96").Value = -10829.333
$ws.Range("N96").Value = -90487
$ws.Range("H99").Value = 4412.5
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 5616.6665
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 5616.6665
$ws.Range("M99").Value = 698
$ws.Range("N99").Value = -8612.666499999999
$ws.Range("H102").Value = 8323.799999999999
$ws.Range("I102").Value = 8323.799999999999
$ws.Range("K102").Value = 8323.799999999999
$ws.Range("M102").Value = -5078.799999999999
$ws.Range("H134").Value = 3443.913
$ws.Range("I134").Value = 3639.6924
$ws.Range("K134").Value = 10919.0772
$ws.Range("M134").Value = -8384.0772

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 32
$ws.Range("K4").Value = 32
$ws.Range("M4").Value = 80
$ws.Range("H50").Value = 46665.832
$ws.Range("J50").Value = 53999
$ws.Range("L50").Value = 53999
$ws.Range("N50").Value = -55249
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472
$ws.Range("H58").Value = 2800.2942
$ws.Range("I58").Value = 2640.4
$ws.Range("K58").Value = 2640.4
$ws.Range("M58").Value = -2437.4
$ws.Range("H60").Value = 32500
$ws.Range("J60").Value = 50000
$ws.Range("L60").Value = 50000
$ws.Range("N60").Value = -51022
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696
$ws.Range("H62").Value = 4159.8
$ws.Range("I62").Value = 4033
$ws.Range("K62").Value = 4033
$ws.Range("M62").Value = -3409
$ws.Range("H65").Value = 4159.8
$ws.Range("I65").Value = 4033
$ws.Range("K65").Value = 20165
$ws.Range("M65").Value = -17045
$ws.Range("H68").Value = 66000
$ws.Range("J68").Value = 67500
$ws.Range("L68").Value = 67500
$ws.Range("N68").Value = -68998
$ws.Range("H71").Value = 66000
$ws.Range("J71").Value = 67500
$ws.Range("L71").Value = 202500
$ws.Range("N71").Value = -209988
$ws.Range("H93").Value = 24299.889
$ws.Range("I93").Value = 23587.375
$ws.Range("K93").Value = 23587.375
$ws.Range("M93").Value = -21715.375
$ws.Range("H111").Value = 32150.5
$ws.Range("J111").Value = 32150.5
$ws.Range("L111").Value = 32150.5
$ws.Range("N111").Value = -40330.5
$ws.Range("H132").Value = 2354.8
$ws.Range("I132").Value = 2354.8
$ws.Range("K132").Value = 7064.400000000001
$ws.Range("M132").Value = -4534.400000000001
$ws.Range("H134").Value = 2074.6667
$ws.Range("I134").Value = 2090.0667
$ws.Range("K134").Value = 6270.2001
$ws.Range("M134").Value = -3735.2001
$ws.Range("H136").Value = 2800.2942
$ws.Range("I136").Value = 2640.4
$ws.Range("K136").Value = 7921.200000000001
$ws.Range("M136").Value = -5371.200000000001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 173.08333
$ws.Range("I10").Value = 184.27272
$ws.Range("J10").Value = 50
$ws.Range("K10").Value = 552.81816
$ws.Range("L10").Value = 150
$ws.Range("M10").Value = -413.81816
$ws.Range("N10").Value = -428
$ws.Range("H140").Value = 3514.8462
$ws.Range("I140").Value = 2522.2222
$ws.Range("J140").Value = 5748.25
$ws.Range("K140").Value = 7566.6666
$ws.Range("L140").Value = 17244.75
$ws.Range("M140").Value = -2386.6666
$ws.Range("N140").Value = -27604.75

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 422.2857
$ws.Range("I2").Value = 303.22223
$ws.Range("J2").Value = 636.6
$ws.Range("K2").Value = 303.22223
$ws.Range("L2").Value = 636.6
$ws.Range("M2").Value = -190.22223
$ws.Range("N2").Value = -862.6
$ws.Range("H23").Value = 2926.5
$ws.Range("I23").Value = 3470.6667
$ws.Range("J23").Value = 2600
$ws.Range("K23").Value = 3470.6667
$ws.Range("L23").Value = 2600
$ws.Range("M23").Value = -3247.6667
$ws.Range("N23").Value = -3046
$ws.Range("H26").Value = 59999.5
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49720
$ws.Range("H50").Value = 59999.5
$ws.Range("I50").Value = 50000
$ws.Range("K50").Value = 50000
$ws.Range("M50").Value = -49502
$ws.Range("H97").Value = 1311.1765
$ws.Range("I97").Value = 1345.375
$ws.Range("K97").Value = 1345.375
$ws.Range("M97").Value = -849.375
$ws.Range("H102").Value = 7351.5
$ws.Range("I102").Value = 7221.8
$ws.Range("K102").Value = 7221.8
$ws.Range("M102").Value = -5599.8
$ws.Range("H132").Value = 2689.25
$ws.Range("I132").Value = 2314.9312
$ws.Range("K132").Value = 6944.7936
$ws.Range("M132").Value = -4414.7936

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2309.4
$ws.Range("I7").Value = 2421.6667
$ws.Range("K7").Value = 2421.6667
$ws.Range("M7").Value = -2309.6667
$ws.Range("H87").Value = 500020000
$ws.Range("J87").Value = 500020000
$ws.Range("L87").Value = 500020000
$ws.Range("N87").Value = -500022246
$ws.Range("H90").Value = 500020000
$ws.Range("J90").Value = 500020000
$ws.Range("L90").Value = 1500060000
$ws.Range("N90").Value = -1500071232
$ws.Range("H97").Value = 23333.334
$ws.Range("J97").Value = 23333.334
$ws.Range("L97").Value = 23333.334
$ws.Range("N97").Value = -25315.334
$ws.Range("H122").Value = 4534.6665
$ws.Range("I122").Value = 4333.857
$ws.Range("J122").Value = 5237.5
$ws.Range("K122").Value = 13001.571
$ws.Range("L122").Value = 15712.5
$ws.Range("M122").Value = -10551.571
$ws.Range("N122").Value = -20612.5
$ws.Range("H126").Value = 2309.4
$ws.Range("I126").Value = 2421.6667
$ws.Range("K126").Value = 7265.000100000001
$ws.Range("M126").Value = -4795.000100000001
$ws.Range("H136").Value = 24702.084
$ws.Range("I136").Value = 4538.154
$ws.Range("J136").Value = 48532.184
$ws.Range("K136").Value = 13614.462
$ws.Range("L136").Value = 145596.552
$ws.Range("M136").Value = -11064.462
$ws.Range("N136").Value = -150696.552

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 49520
$ws.Range("J94").Value = 49520
$ws.Range("L94").Value = 49520
$ws.Range("N94").Value = -51322
$ws.Range("H95").Value = 20344
$ws.Range("J95").Value = 20344
$ws.Range("L95").Value = 20344
$ws.Range("N95").Value = -25836
$ws.Range("H96").Value = 3984.9333
$ws.Range("I96").Value = 2897.1428
$ws.Range("J96").Value = 4936.75
$ws.Range("K96").Value = 2897.1428
$ws.Range("L96").Value = 4936.75
$ws.Range("M96").Value = -1524.1428
$ws.Range("N96").Value = -7682.75
